$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") contains a date serial value (45175 -> 2023-09-06)
# that was bumped to 45177 (-> 2023-09-08) for every data row (2..480).
$startRow = 2
$endRow = 480

$rng = $ws.Range("C$startRow`:C$endRow")
$rng.Value = 45177
